# Fixed some biomass transcription errors.
#
# 1) The empty paragraph that follows "9 total samples" currently holds the
#    hidden "_GoBack" bookmark left over from the last edit location. It
#    should just become a plain empty paragraph.
# 2) A new trailing paragraph of bold/italic text describing the sampling
#    schedule is added at the very end of the document (after "12 total
#    samples"), preceded by a blank paragraph, and the "_GoBack" bookmark
#    now marks that new final location instead.

$d = $word.ActiveDocument

# --- Step 1: strip the _GoBack bookmark out of the paragraph after
#     "9 total samples", leaving a bare empty paragraph in its place.
$goBackPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq "9 total samples`r") {
        $goBackPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

$emptyParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p/></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$goBackPara.Range.InsertXML($emptyParagraphXml)

# --- Step 2: append a blank paragraph plus the new bold/italic sampling
#     note (carrying the _GoBack bookmark) at the very end of the document.
$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos, $endPos)

$newTailXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p/>' +
    '<w:p>' +
    '<w:pPr><w:rPr><w:b/><w:i/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>Sampling in June, September, November, January</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newTailXml)

Write-Output "Applied biomass transcription fixes."
